$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 247.3125
$ws.Range("I33").Value = 239.85715
$ws.Range("J33").Value = 299.5
$ws.Range("K33").Value = 239.85715
$ws.Range("L33").Value = 299.5
$ws.Range("M33").Value = -10.85714999999999
$ws.Range("N33").Value = -757.5
$ws.Range("H40").Value = 5317.2354
$ws.Range("I40").Value = 9499.166999999999
$ws.Range("K40").Value = 9499.166999999999
$ws.Range("M40").Value = -9324.166999999999
$ws.Range("H125").Value = 4764
$ws.Range("I125").Value = 2173.75
$ws.Range("J125").Value = 8217.666999999999
$ws.Range("K125").Value = 19563.75
$ws.Range("L125").Value = 73959.003
$ws.Range("M125").Value = -17103.75
$ws.Range("N125").Value = -78879.003
$ws.Range("H129").Value = 2018.25
$ws.Range("I129").Value = 1806.5714
$ws.Range("K129").Value = 5419.7142
$ws.Range("M129").Value = -419.7142000000003
$ws.Range("H132").Value = 49710.715
$ws.Range("I132").Value = 56991.453
$ws.Range("J132").Value = 6026.2856
$ws.Range("K132").Value = 170974.359
$ws.Range("L132").Value = 18078.8568
$ws.Range("M132").Value = -168444.359
$ws.Range("N132").Value = -23138.8568
$ws.Range("H137").Value = 2320622.8
$ws.Range("I137").Value = 1468
$ws.Range("K137").Value = 4404
$ws.Range("M137").Value = -1854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5095.9443
$ws.Range("I2").Value = 3706.1428
$ws.Range("K2").Value = 3706.1428
$ws.Range("M2").Value = -3593.1428
$ws.Range("H32").Value = 11906273
$ws.Range("I32").Value = 13159370
$ws.Range("K32").Value = 13159370
$ws.Range("M32").Value = -13159083
$ws.Range("H45").Value = 3135.2273
$ws.Range("I45").Value = 3127.0557
$ws.Range("J45").Value = 3172
$ws.Range("K45").Value = 3127.0557
$ws.Range("L45").Value = 3172
$ws.Range("M45").Value = -2750.0557
$ws.Range("N45").Value = -3926
$ws.Range("H63").Value = 20977.834
$ws.Range("I63").Value = 5172
$ws.Range("K63").Value = 5172
$ws.Range("M63").Value = -4486
$ws.Range("H66").Value = 20977.834
$ws.Range("I66").Value = 5172
$ws.Range("K66").Value = 25860
$ws.Range("M66").Value = -22428
$ws.Range("H116").Value = 5095.9443
$ws.Range("I116").Value = 3706.1428
$ws.Range("K116").Value = 3706.1428
$ws.Range("M116").Value = -1412.1428
$ws.Range("H122").Value = 2486.4666
$ws.Range("I122").Value = 1433.1111
$ws.Range("J122").Value = 4066.5
$ws.Range("K122").Value = 4299.3333
$ws.Range("L122").Value = 12199.5
$ws.Range("M122").Value = -1849.3333
$ws.Range("N122").Value = -17099.5
$ws.Range("H124").Value = 40214.5
$ws.Range("J124").Value = 40214.5
$ws.Range("L124").Value = 40214.5
$ws.Range("N124").Value = -50034.5
$ws.Range("H125").Value = 64426.25
$ws.Range("J125").Value = 64426.25
$ws.Range("L125").Value = 64426.25
$ws.Range("N125").Value = -74266.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5095.9443
$ws.Range("I3").Value = 3706.1428
$ws.Range("K3").Value = 3706.1428
$ws.Range("M3").Value = -3592.1428
$ws.Range("H20").Value = 1823.5555
$ws.Range("I20").Value = 1494.6666
$ws.Range("K20").Value = 1494.6666
$ws.Range("M20").Value = -1247.6666
$ws.Range("H43").Value = 684000
$ws.Range("J43").Value = 684000
$ws.Range("L43").Value = 684000
$ws.Range("N43").Value = -684362
$ws.Range("H86").Value = 2060.7856
$ws.Range("I86").Value = 1875.7
$ws.Range("J86").Value = 2523.5
$ws.Range("K86").Value = 1875.7
$ws.Range("L86").Value = 2523.5
$ws.Range("M86").Value = -752.7
$ws.Range("N86").Value = -4769.5
$ws.Range("H89").Value = 2060.7856
$ws.Range("I89").Value = 1875.7
$ws.Range("J89").Value = 2523.5
$ws.Range("K89").Value = 9378.5
$ws.Range("L89").Value = 12617.5
$ws.Range("M89").Value = -3762.5
$ws.Range("N89").Value = -23849.5
$ws.Range("H94").Value = 818.7143
$ws.Range("I94").Value = 685.5
$ws.Range("J94").Value = 1058.5
$ws.Range("K94").Value = 685.5
$ws.Range("L94").Value = 1058.5
$ws.Range("M94").Value = -234.5
$ws.Range("N94").Value = -1960.5
$ws.Range("H99").Value = 7381.391
$ws.Range("I99").Value = 3680.353
$ws.Range("K99").Value = 3680.353
$ws.Range("M99").Value = -2182.353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 876.8570999999999
$ws.Range("I22").Value = 856.3333
$ws.Range("K22").Value = 856.3333
$ws.Range("M22").Value = -506.3333
$ws.Range("H58").Value = 953633.6
$ws.Range("I58").Value = 1373193.1
$ws.Range("K58").Value = 1373193.1
$ws.Range("M58").Value = -1372990.1
$ws.Range("H86").Value = 124607.53
$ws.Range("I86").Value = 6767.091
$ws.Range("J86").Value = 340648.34
$ws.Range("K86").Value = 6767.091
$ws.Range("L86").Value = 340648.34
$ws.Range("M86").Value = -5644.091
$ws.Range("N86").Value = -342894.34
$ws.Range("H89").Value = 124607.53
$ws.Range("I89").Value = 6767.091
$ws.Range("J89").Value = 340648.34
$ws.Range("K89").Value = 33835.455
$ws.Range("L89").Value = 1703241.7
$ws.Range("M89").Value = -28219.455
$ws.Range("N89").Value = -1714473.7
$ws.Range("H93").Value = 36553.89
$ws.Range("J93").Value = 42500
$ws.Range("L93").Value = 42500
$ws.Range("N93").Value = -46244
$ws.Range("H124").Value = 41796.934
$ws.Range("J124").Value = 41796.934
$ws.Range("L124").Value = 41796.934
$ws.Range("N124").Value = -46706.934
$ws.Range("H136").Value = 953633.6
$ws.Range("I136").Value = 1373193.1
$ws.Range("K136").Value = 4119579.3
$ws.Range("M136").Value = -4117029.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6499.857
$ws.Range("I131").Value = 626.2857
$ws.Range("J131").Value = 7674.5713
$ws.Range("K131").Value = 1878.8571
$ws.Range("L131").Value = 23023.7139
$ws.Range("M131").Value = 3161.1429
$ws.Range("N131").Value = -33103.7139
$ws.Range("H134").Value = 632.0833
$ws.Range("I134").Value = 632.0833
$ws.Range("K134").Value = 1896.2499
$ws.Range("M134").Value = 3173.7501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5214
$ws.Range("I70").Value = 4999.75
$ws.Range("J70").Value = 5499.6665
$ws.Range("K70").Value = 4999.75
$ws.Range("L70").Value = 5499.6665
$ws.Range("M70").Value = -4729.75
$ws.Range("N70").Value = -6039.6665
$ws.Range("H73").Value = 5214
$ws.Range("I73").Value = 4999.75
$ws.Range("J73").Value = 5499.6665
$ws.Range("K73").Value = 4999.75
$ws.Range("L73").Value = 5499.6665
$ws.Range("M73").Value = -4063.75
$ws.Range("N73").Value = -7371.6665
$ws.Range("H122").Value = 123590.78
$ws.Range("I122").Value = 254328
$ws.Range("K122").Value = 762984
$ws.Range("M122").Value = -760534

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 74059.57000000001
$ws.Range("J63").Value = 75068
$ws.Range("L63").Value = 75068
$ws.Range("N63").Value = -76566
$ws.Range("H66").Value = 74059.57000000001
$ws.Range("J66").Value = 75068
$ws.Range("L66").Value = 225204
$ws.Range("N66").Value = -232692
$ws.Range("H69").Value = 70163
$ws.Range("J69").Value = 70163
$ws.Range("L69").Value = 70163
$ws.Range("N69").Value = -71785
$ws.Range("H72").Value = 70163
$ws.Range("J72").Value = 70163
$ws.Range("L72").Value = 210489
$ws.Range("N72").Value = -218601

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3175.6562
$ws.Range("I122").Value = 2825.9285
$ws.Range("J122").Value = 5623.75
$ws.Range("K122").Value = 8477.7855
$ws.Range("L122").Value = 16871.25
$ws.Range("M122").Value = -6027.7855
$ws.Range("N122").Value = -21771.25
$ws.Range("H136").Value = 1592746.8
$ws.Range("I136").Value = 2038916.9
$ws.Range("K136").Value = 6116750.699999999
$ws.Range("M136").Value = -6114200.699999999
